# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price (col D) and 1h volume-change (col E) values are refreshed per coin row;
# rows 45-47 additionally get reshuffled coin identities (name/link) because the
# underlying ranking order changed between scrapes.
#
# Values are forced to Text (NumberFormat "@") before being written so numeric-
# looking strings such as "1.00" or "302.40" are not silently reinterpreted as
# numbers by Excel; the style is then reset to "Normal" so no stray formatting
# is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '43.066.23'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.11%  '

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.314.83'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +0.04%  '

$ws.Range('E4').Value = '  +0.00%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '302.40'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.29%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '99.20'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -1.91%  '

$ws.Range('E7').Value = '  +0.28%  '

$ws.Range('E8').Value = '  +0.00%  '

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.521'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +2.12%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '36.01'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +3.00%  '

$ws.Range('E11').Value = '  -0.96%  '

$ws.Range('E12').Value = '  -1.28%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '17.76'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -1.30%  '

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '6.88'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +0.18%  '

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '2.675.82'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -0.58%  '

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '2.299.30'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -2.37%  '

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.790'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -3.80%  '

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '42.995.96'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.14%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '13.20'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +5.53%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '6.19'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.21%  '

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.0₃0908'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.17%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '68.16'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.38%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '240.93'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +1.46%  '

$ws.Range('E24').Value = '  -3.57%  '

$ws.Range('E25').Value = '  -1.09%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.07%  '

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '25.15'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +1.19%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '169.53'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +0.70%  '

$ws.Range('E29').Value = '  -2.08%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '9.17'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.35%  '

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '33.32'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -2.76%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '4.97'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +6.85%  '

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '5.19'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +2.88%  '

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.15%  '

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '18.37'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +6.67%  '

$ws.Range('E36').Value = '  -0.81%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.0694'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -0.23%  '

$ws.Range('E38').Value = '  +0.95%  '

$ws.Range('E39').Value = '  -0.37%  '

$ws.Range('E40').Value = '  -2.50%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.110'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.36%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.996.47'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -0.46%  '

$ws.Range('E43').Value = '  +0.43%  '

$ws.Range('E44').Value = '  -0.77%  '

$c = $ws.Range('B45')
$c.NumberFormat = '@'
$c.Value = 'ApeXProtocol'
$c.Style = 'Normal'
$c = $ws.Range('C45')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c.Style = 'Normal'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.09'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -9.68%  '

$c = $ws.Range('B46')
$c.NumberFormat = '@'
$c.Value = 'EnergySwap'
$c.Style = 'Normal'
$c = $ws.Range('C46')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '17.42'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -1.47%  '

$c = $ws.Range('B47')
$c.NumberFormat = '@'
$c.Value = 'NEARProtocol'
$c.Style = 'Normal'
$c = $ws.Range('C47')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '2.84'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -1.03%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '76.28'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +8.27%  '

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '54.91'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -2.16%  '

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.542.82'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +0.62%  '

$ws.Range('E51').Value = '  -0.68%  '
